# Update automatico via Actualizar 09-23-2020 00-42-22
# Adds a new "Teletrabajo" source row (row 42) to the DATACOVID "trabajo" table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Grow the table by one row - this extends the table ref, autoFilter ref,
# and keeps everything wired up to the ListObject.
$lo.ListRows.Add() | Out-Null

# Copy the formatting (styles + row height) from the last existing data
# row (41) down onto the freshly added row (42) so the new row matches
# the look of the rest of the table.
$ws.Range("A41:K41").Copy()
$ws.Range("A42:K42").PasteSpecial(-4122)
$ws.Rows.Item(42).RowHeight = $ws.Rows.Item(41).RowHeight

# Columns that repeat the same source/category/country info as row 41 -
# read back via .Text (not .Value) so the exact existing shared-string
# text (incl. any non-breaking spaces) is reused instead of minted anew.
$ws.Range("A42").Value = $ws.Range("A41").Text
$ws.Range("C42").Value = $ws.Range("C41").Text
$ws.Range("D42").Value = $ws.Range("D41").Text
$ws.Range("G42").Value = $ws.Range("G41").Text
$ws.Range("J42").Value = $ws.Range("J41").Text
$ws.Range("K42").Value = $ws.Range("K41").Text

# Brand-new values for this row.
$ws.Range("B42").Value = 41
$ws.Range("E42").Value = "https://www.mitradel.gob.pa/decreto-ejecutivo-que-reglamenta-la-ley-de-teletrabajo-en-gaceta-oficial/"
$ws.Range("F42").Value = "El Decreto Ejecutivo Número 133 del 16 de septiembre de 2020, que reglamenta la Ley de Teletrabajo fue publicado en Gaceta Oficial luego que el presidente de la República, Laurentino Cortizo Cohen, firmara el documento en un acto que oficializaba todos los acuerdos alcanzados en la Mesa Tripartita de Diálogo por la Economía y el Desarrollo Laboral."
$ws.Range("H42").Value = 44095
$ws.Range("I42").Value = 44095

# Wire up the hyperlinks for the "Sitio Web" and "Descarga Link" columns,
# then re-apply the table formatting on top (Hyperlinks.Add swaps in the
# built-in "Hyperlink" style, which we don't want here).
$ws.Hyperlinks.Add($ws.Range("G42"), "https://www.mitradel.gob.pa/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E42"), "https://www.mitradel.gob.pa/decreto-ejecutivo-que-reglamenta-la-ley-de-teletrabajo-en-gaceta-oficial/") | Out-Null

$ws.Range("G41").Copy()
$ws.Range("G42").PasteSpecial(-4122)
$ws.Range("E41").Copy()
$ws.Range("E42").PasteSpecial(-4122)

# Extend the category data-validation down to the new row.
$ws.Range("C2:C41").Validation.Delete()
$v = $ws.Range("C2:C42").Validation
$v.Add(0)
$v.ErrorTitle = "Entrada no válida"
$v.ErrorMessage = "Selecciona una categoría de la lista"
$v.InputTitle = "Categoria"
$v.InputMessage = "Selecciona una categoría de la lista"

# Match the author's final selection (cell I42).
$ws.Range("I42").Select() | Out-Null
